$wb = $excel.ActiveWorkbook

# --- Sheet "Test Cases" ---
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("C2").Value = "N"
$ws1.Range("D2").Value = "PASS"
$ws1.Range("C3").Value = "Y"
$ws1.Range("D3").Value = "PASS"

# --- Sheet "Checklogin" ---
$ws2 = $wb.Worksheets.Item("Checklogin")
$ws2.Range("C3").Value = "N"
$ws2.Range("D3").Value = "SKIP"
$ws2.Range("C4").Value = "N"
$ws2.Range("D4").Value = "SKIP"
$ws2.Range("C5").Value = "N"
$ws2.Range("D5").Value = "SKIP"

# --- Sheet "NewCustomerRegistration" ---
$ws3 = $wb.Worksheets.Item("NewCustomerRegistration")
# Drop the RadioButton..PreviousOptician columns (old H:N), then PreferredName (old D)
$ws3.Range("H1:N2").EntireColumn.Delete()
$ws3.Range("D1:D2").EntireColumn.Delete()

$ws3.Range("A2").Value = "MR"
$ws3.Range("B2").Value = "Gupte"
$ws3.Range("C2").Value = "Aakar"
$ws3.Range("D2").Value = "26/04/1990"
$ws3.Range("E2").Value = "M"
$ws3.Range("F2").Value = 1236547890
$ws3.Range("G2").Value = "Y"
$ws3.Range("H2").Value = "PASS"

# --- Selections / active sheet ---
$ws1.Range("D2").Select()
$ws2.Range("C6").Select()
$ws3.Range("A3").Select()

$ws1.Activate()
